# ncp-gop-transect-info.xlsx -- metadata updates and handling incorrect datetimes
#
# 1. Personnel sheet: Rachel Stanley's role is corrected from
#    "principal Investigator" to the controlled-vocabulary value "PI".
# 2. Keywords sheet: two new controlled keywords are appended
#    ("net ecosystem production", "net community production"), matching
#    the wrapped/centered formatting already used for the other long
#    keyword rows.
# 3. The active worksheet/selection moves from Personnel to Keywords.

$wb  = $excel.ActiveWorkbook
$personnel = $wb.Worksheets.Item("Personnel")
$keywords  = $wb.Worksheets.Item("Keywords")

# --- Personnel: role correction --------------------------------------
$personnel.Range("G5").Value = "PI"

# --- Keywords: append two new keyword rows ----------------------------
$keywords.Range("A9").Value = "net ecosystem production"
$keywords.Range("A9").WrapText = $true
$keywords.Range("A9").VerticalAlignment = -4108
$keywords.Range("A9").RowHeight = 31.2

$keywords.Range("A10").Value = "net community production"
$keywords.Range("A10").WrapText = $true
$keywords.Range("A10").VerticalAlignment = -4108
$keywords.Range("A10").RowHeight = 31.2

# --- Selection / active sheet bookkeeping ------------------------------
[void]$personnel.Range("G12").Select()

[void]$keywords.Activate()
[void]$keywords.Range("A9:A10").Select()
